# Execution.xlsx - "Add files via upload"
# Update the Summary sheet's application list: rename the GL/CX rows to their
# correct master-sheet names, and add five new application rows (Manufacturing
# was inserted where CX used to be, then Inventory / Costing / My Receiving /
# Order Management are appended).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- Row 6: GL -------------------------------------------------------------
$ws.Range("E6").Value = "Master_GL.xlsx"
$ws.Range("H6").Style = "Normal"
$ws.Range("H6").Value = "NO"

# --- Row 7: CX -> Manufacturing --------------------------------------------
$ws.Range("C7").Value = "Manufacturing"
$ws.Range("E7").Value = "Master_Manufacturing.xlsx"
$ws.Range("F7").Style = "Normal"
$ws.Range("F7").Value = "Reusable_Components.xlsx"
$ws.Range("G7").Style = "Normal"
$ws.Range("G7").Value = "Reusable_Components"
$ws.Range("H7").Style = "Normal"
$ws.Range("H7").Value = "NO"

# --- Row 8: AP (values stay the same, formatting + blank cells change) ----
$ws.Range("F8").Value = ""
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = ""
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Style = "Normal"
$ws.Range("H8").Value = "NO"

# --- New rows 9-12 -----------------------------------------------------
$ws.Range("A9:H12").Style = "Normal"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Selenium"
$ws.Range("C9").Value = "Inventory"
$ws.Range("D9").Value = "Web"
$ws.Range("E9").Value = "Master_Inventory.xlsx"
$ws.Range("F9").Value = "Reusable_Components.xlsx"
$ws.Range("G9").Value = "Reusable_Components"
$ws.Range("H9").Value = "NO"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Selenium"
$ws.Range("C10").Value = "Costing"
$ws.Range("D10").Value = "Web"
$ws.Range("E10").Value = "Master_Costing.xlsx"
$ws.Range("F10").Value = "Reusable_Components.xlsx"
$ws.Range("G10").Value = "Reusable_Components"
$ws.Range("H10").Value = "YES"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Selenium"
$ws.Range("C11").Value = "My Receiving"
$ws.Range("D11").Value = "Web"
$ws.Range("E11").Value = "Master_MyReceiving_20D.xlsx"
$ws.Range("F11").Value = "Reusable_Components.xlsx"
$ws.Range("G11").Value = "Reusable_Components"
$ws.Range("H11").Value = "NO"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Selenium"
$ws.Range("C12").Value = "Order Management"
$ws.Range("D12").Value = "Web"
$ws.Range("E12").Value = "Master_OrderManagement.xlsx"
$ws.Range("F12").Value = "Reusable_Components.xlsx"
$ws.Range("G12").Value = "Reusable_Components"
$ws.Range("H12").Value = "NO"

# --- Selection ends on H12 (last edited cell) -------------------------
$ws.Range("H12").Select()
